$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.248976349830627
$ws.Range("B1").Value = 2.214832544326782
$ws.Range("C1").Value = 2.825379848480225
$ws.Range("D1").Value = 3.275952577590942
$ws.Range("E1").Value = 2.118652582168579
